# Re-apply the latest "cryptos list" scrape values (GitHub Actions bot commit).
# Only the Price (column D) and Volume(1h) (column E) cells for specific rows
# change; everything else (labels, links, rank index, formatting) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several of the new Price values look like plain numbers (e.g. "258.72",
# "1.000"). Excel's Range.Value setter auto-converts such strings to numeric
# values, which would turn "1.000" into 1 and drop the trailing zeros. The
# source data keeps these as literal text, so force the whole Price/Volume
# block to Text format before writing, then restore normal styling
# afterwards so no stray per-cell format sticks around.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '26.581.85'
$ws.Range("D3").Value = '1.839.57'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '258.72'
$ws.Range("E5").Value = '  -1.27%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '0.5276'
$ws.Range("E7").Value = '  +1.09%  '
$ws.Range("D8").Value = '0.3139'
$ws.Range("E8").Value = '  -3.65%  '
$ws.Range("D9").Value = '0.06795'
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("D10").Value = '18.66'
$ws.Range("E10").Value = '  -0.15%  '
$ws.Range("D11").Value = '0.7782'
$ws.Range("E11").Value = '  +0.75%  '
$ws.Range("D12").Value = '0.07756'
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").Value = '1.838.66'
$ws.Range("E13").Value = '  -0.21%  '
$ws.Range("D14").Value = '87.75'
$ws.Range("E14").Value = '  -0.31%  '
$ws.Range("D15").Value = '5.001'
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("E17").Value = '  -0.74%  '
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").Value = '0.000007902'
$ws.Range("E19").Value = '  -0.67%  '
$ws.Range("D20").Value = '26.611.17'
$ws.Range("D21").Value = '2.081.66'
$ws.Range("E21").Value = '  +0.59%  '
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").Value = '5.960'
$ws.Range("E23").Value = '  -0.59%  '
$ws.Range("D24").Value = '9.303'
$ws.Range("E24").Value = '  -2.33%  '
$ws.Range("D25").Value = '142.94'
$ws.Range("E25").Value = '  -1.31%  '
$ws.Range("D26").Value = '2.203'
$ws.Range("E26").Value = '  +1.22%  '
$ws.Range("D27").Value = '1.681'
$ws.Range("E27").Value = '  +1.61%  '
$ws.Range("D29").Value = '110.47'
$ws.Range("E29").Value = '  -1.05%  '
$ws.Range("D30").Value = '4.173'
$ws.Range("E30").Value = '  -0.55%  '
$ws.Range("D31").Value = '0.08719'
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("D32").Value = '4.059'
$ws.Range("E32").Value = '  -1.69%  '
$ws.Range("D33").Value = '0.04870'
$ws.Range("D34").Value = '0.7297'
$ws.Range("E34").Value = '  +1.60%  '
$ws.Range("E35").Value = '  +0.50%  '
$ws.Range("D36").Value = '2.859'
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("D37").Value = '3.090'
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("D38").Value = '2.244'
$ws.Range("E38").Value = '  +0.76%  '
$ws.Range("D39").Value = '0.01720'
$ws.Range("E39").Value = '  -3.37%  '
$ws.Range("D40").Value = '0.4785'
$ws.Range("E40").Value = '  -1.19%  '
$ws.Range("D41").Value = '0.8941'
$ws.Range("E41").Value = '  -0.74%  '
$ws.Range("D42").Value = '109.89'
$ws.Range("E42").Value = '  -2.21%  '
$ws.Range("D43").Value = '5.912'
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").Value = '7.643'
$ws.Range("E45").Value = '  -1.12%  '
$ws.Range("D46").Value = '0.4150'
$ws.Range("E46").Value = '  +0.17%  '
$ws.Range("D47").Value = '8.969'
$ws.Range("E47").Value = '  +0.18%  '
$ws.Range("D48").Value = '0.1235'
$ws.Range("E48").Value = '  +0.71%  '
$ws.Range("E49").Value = '  -1.25%  '
$ws.Range("D50").Value = '34.66'
$ws.Range("E50").Value = '  -1.08%  '
$ws.Range("D51").Value = '0.8938'
$ws.Range("E51").Value = '  +0.87%  '

# Restore default (Normal) styling so no stray per-cell style index is left behind.
$dataRange.Style = "Normal"
